# Excel constants (named, for readability - values are the standard COM enum values)
$xlCenter      = -4108   # XlHAlign / XlVAlign
$xlTop         = -4160   # XlVAlign
$xlContinuous  = 1       # XlLineStyle
$xlPasteFormats = -4122  # XlPasteType

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (Plan1 -> Sheet1)
$ws.Name = "Sheet1"

# Populate the keyword data (header + 3 rows of numbered keywords)
$ws.Range("B1").Value = "times"
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "flu"
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "vasco"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "fla"

# Format the header/index cells: bold, centered horizontally, top-aligned vertically,
# with a thin box border all around.
$r = $ws.Range("B1")
$r.Font.Bold = $true
$r.HorizontalAlignment = $xlCenter
$r.VerticalAlignment = $xlTop
$r.Borders.LineStyle = $xlContinuous

# Copy that formatting onto A2:A4 so every cell shares the same style
# (copy/paste-special reuses the style instead of registering new ones per cell).
$r.Copy()
$ws.Range("A2").PasteSpecial($xlPasteFormats)
$ws.Range("A3").PasteSpecial($xlPasteFormats)
$ws.Range("A4").PasteSpecial($xlPasteFormats)
